$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.551.01"
$ws.Range("E2").Value = "  +1.13%  "

# Row 3
$ws.Range("D3").Value = "1.761.60"
$ws.Range("E3").Value = "  -1.34%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.87"
$ws.Range("E5").Value = "  +0.26%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3833"
$ws.Range("E7").Value = "  +0.70%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3402"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.92"
$ws.Range("E9").Value = "  -2.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.137"
$ws.Range("E10").Value = "  -4.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07394"
$ws.Range("E11").Value = "  -1.11%  "

# Row 12
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.33"
$ws.Range("E13").Value = "  +2.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.339"
$ws.Range("E14").Value = "  -1.79%  "

# Row 15
$ws.Range("D15").Value = "1.761.62"
$ws.Range("E15").Value = "  -1.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.026"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17
$ws.Range("E17").Value = "  -2.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06656"
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.18"
$ws.Range("E19").Value = "  -1.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.32"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.388"
$ws.Range("E22").Value = "  -3.55%  "

# Row 23
$ws.Range("D23").Value = "27.559.28"
$ws.Range("E23").Value = "  +1.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.02"
$ws.Range("E24").Value = "  -2.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.376"
$ws.Range("E25").Value = "  -1.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.428"
$ws.Range("E26").Value = "  -4.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  -3.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.411"
$ws.Range("E28").Value = "  -4.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.51"
$ws.Range("E29").Value = "  -0.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.34"
$ws.Range("E30").Value = "  +0.24%  "

# Row 31
$ws.Range("D31").Value = "1.963.15"
$ws.Range("E31").Value = "  -1.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.109"
$ws.Range("E32").Value = "  +0.74%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.959"
$ws.Range("E33").Value = "  -1.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08787"
$ws.Range("E34").Value = "  +0.93%  "

# Row 35
$ws.Range("E35").Value = "  -4.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02409"
$ws.Range("E36").Value = "  +3.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6767"
$ws.Range("E37").Value = "  -2.40%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.325"
$ws.Range("E38").Value = "  -2.04%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06289"
$ws.Range("E39").Value = "  -0.53%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2179"
$ws.Range("E40").Value = "  -1.26%  "

# Row 41
$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.510"
$ws.Range("E41").Value = "  -8.83%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.243"
$ws.Range("E42").Value = "  +0.57%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.255"

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.20"
$ws.Range("E44").Value = "  -1.64%  "

# Row 45
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6244"
$ws.Range("E46").Value = "  -4.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.823"
$ws.Range("E47").Value = "  -0.65%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.28"
$ws.Range("E48").Value = "  +1.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.074"
$ws.Range("E49").Value = "  -3.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07367"
$ws.Range("E50").Value = "  +3.34%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.145"
$ws.Range("E51").Value = "  +2.43%  "
